# feat: add 2022-Q3 data
#
# The "2021-Q1" sheet's data becomes the new "2022-Q3" quarter (renamed +
# overwritten with the new quarter's numbers), and a fresh copy of the
# original "2021-Q1" content is kept (under its original name) right after
# it. The summary ("总计") sheet gets a new row for 2022-Q3 inserted ahead
# of the existing 2021-Q1 total row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. 总计 (summary) sheet: insert a 2022-Q3 row before the 2021-Q1 row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Make room below the existing data row (row 2) for the old 2021-Q1 values,
# which move down to row 3.
$summary.Rows.Item(3).Insert()

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = $summary.Cells.Item(2,2).Value2
$summary.Cells.Item(3,3).Value = $summary.Cells.Item(2,3).Value2
$summary.Cells.Item(3,4).Value = $summary.Cells.Item(2,4).Value2

# Match the header/index cell formatting (bold, centered, thin border) that
# row 2 already carries.
$summary.Cells.Item(3,1).Font.Bold = $true
$summary.Cells.Item(3,1).HorizontalAlignment = -4108
$summary.Cells.Item(3,1).VerticalAlignment = -4160
$summary.Cells.Item(3,1).Borders.LineStyle = 1

# Row 2 now becomes the new 2022-Q3 total row.
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 1
$summary.Cells.Item(2,4).Value = 0.13

# ---------------------------------------------------------------------
# 2. Quarter detail sheets: duplicate the current "2021-Q1" sheet so the
#    original data is preserved under its own tab, then turn the original
#    tab into "2022-Q3" with the new quarter's fund data.
# ---------------------------------------------------------------------
$oldQ1 = $wb.Worksheets.Item(2)

# Copy placed right after $oldQ1 -- this keeps the original 2021-Q1 content
# byte-for-byte (including formatting) under a new tab.
$oldQ1.Copy($null, $oldQ1)

$oldQ1.Name = "2022-Q3"
$newQ1 = $wb.Worksheets.Item(3)
$newQ1.Name = "2021-Q1"

# $oldQ1 (now named "2022-Q3") had 3 data rows (rows 2-4); the new quarter
# only has a single fund, so drop the extra rows.
$oldQ1.Rows.Item(3).Delete()
$oldQ1.Rows.Item(3).Delete()

$oldQ1.Cells.Item(1,4).Value = "基金规模"

$oldQ1.Cells.Item(2,1).Value = 0
$oldQ1.Cells.Item(2,2).Value = "'004685"
$oldQ1.Cells.Item(2,3).Value = "金元顺安元启灵活配置混合"
$oldQ1.Cells.Item(2,4).Value = "'15.28"
$oldQ1.Cells.Item(2,5).Value = "'77.14"
$oldQ1.Cells.Item(2,6).Value = "'0.86"
$oldQ1.Cells.Item(2,7).Value = "'0.1314"
$oldQ1.Cells.Item(2,8).Value = 8
